$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "halosalsa3" test data rows to "halosalsa4"
$ws.Range("B7").Value = "halosalsa4"
$ws.Range("B8").Value = "halosalsa4@gmail.com"
